$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates from the cryptos list refresh.
# Cells whose new text parses as a plain number get NumberFormat "@"
# applied first (so Excel keeps them as text instead of coercing to a
# number), then ClearFormats() afterwards to drop the formatting again
# so the cell style matches the original (text value is unaffected).

$ws.Range("D2").Value = "70.549.45"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "3.511.83"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "623.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "3.509.12"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.586"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "4.080.21"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.39"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "608.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "3.508.36"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").Value = "70.649.07"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.880"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.11"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.55"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.26"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.03"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.81"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "623.76"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0491"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.80"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0994"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.19%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.65"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.143"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "3.340.85"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "0.0₃0721"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.311"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("E48").Value = "  -5.06%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.22"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.17%  "
